$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mapping")

$ws.Range("A2").Value = -78.9686
$ws.Range("B2").Value = -78.3756

$ws.Range("A3").Value = 33.4733
$ws.Range("B3").Value = 33.9687

$ws.Range("A4").Value = -74.7919
$ws.Range("B4").Value = -75.4151

$ws.Range("A5").Value = 36.8443
$ws.Range("B5").Value = 36.3582
